$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 95-123 (value corrections per diff) ---
# Row 95
$ws.Range("D95").Value = 44559
$ws.Range("J95").Value = 3000
$ws.Range("K95").Value = 2300
$ws.Range("L95").Value = 2300
$ws.Range("M95").Value = 2300
$ws.Range("P95").Value = 2300

# Row 96
$ws.Range("D96").Value = 44559
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = 2000
$ws.Range("P96").Value = 2000

# Row 97
$ws.Range("D97").Value = 44559
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 1500
$ws.Range("P97").Value = 1500

# Row 98
$ws.Range("D98").Value = 44208
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 2000
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 2000
$ws.Range("P98").Value = 2000

# Row 99
$ws.Range("D99").Value = 44208
$ws.Range("J99").Value = 5000

# Row 100
$ws.Range("D100").Value = 44208
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = 1200
$ws.Range("P100").Value = 1200

# Row 101
$ws.Range("D101").Value = 44264
$ws.Range("K101").Value = 2300
$ws.Range("L101").Value = 2300
$ws.Range("M101").Value = 2300
$ws.Range("P101").Value = 2300

# Row 102
$ws.Range("D102").Value = 44264
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1800
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = 1800
$ws.Range("P102").Value = 1800

# Row 103
$ws.Range("D103").Value = 44264
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 1300
$ws.Range("L103").Value = 1300
$ws.Range("M103").Value = 1300
$ws.Range("P103").Value = 1300

# Row 104
$ws.Range("D104").Value = 44232
$ws.Range("J104").Value = 2000
$ws.Range("K104").Value = 2500
$ws.Range("L104").Value = 2500
$ws.Range("M104").Value = 2500
$ws.Range("P104").Value = 2500

# Row 105
$ws.Range("D105").Value = 44232
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 2000
$ws.Range("P105").Value = 2000

# Row 106
$ws.Range("D106").Value = 44232
$ws.Range("I106").Value = 'Segunda'
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1600
$ws.Range("L106").Value = 1600
$ws.Range("M106").Value = 1600
$ws.Range("P106").Value = 1600

# Row 107
$ws.Range("D107").Value = 44551
$ws.Range("I107").Value = 'Extra'
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 3000
$ws.Range("P107").Value = 3000

# Row 108
$ws.Range("D108").Value = 44551
$ws.Range("I108").Value = 'Primera'
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 2500
$ws.Range("L108").Value = 2500
$ws.Range("M108").Value = 2500
$ws.Range("P108").Value = 2500

# Row 109
$ws.Range("D109").Value = 44196
$ws.Range("J109").Value = 1500
$ws.Range("K109").Value = 3000
$ws.Range("L109").Value = 3000
$ws.Range("M109").Value = 3000
$ws.Range("P109").Value = 3000

# Row 110
$ws.Range("D110").Value = 44196
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 2500
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = 2500
$ws.Range("P110").Value = 2500

# Row 111
$ws.Range("D111").Value = 44196
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 2000
$ws.Range("L111").Value = 2000
$ws.Range("M111").Value = 2000
$ws.Range("P111").Value = 2000

# Row 112
$ws.Range("D112").Value = 44200
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 2200
$ws.Range("L112").Value = 2200
$ws.Range("M112").Value = 2200
$ws.Range("P112").Value = 2200

# Row 113
$ws.Range("D113").Value = 44200
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 1700
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = 1700
$ws.Range("P113").Value = 1700

# Row 114
$ws.Range("D114").Value = 44200
$ws.Range("J114").Value = 4000
$ws.Range("K114").Value = 1400
$ws.Range("L114").Value = 1400
$ws.Range("M114").Value = 1400
$ws.Range("P114").Value = 1400

# Row 115
$ws.Range("D115").Value = 44188
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 3500
$ws.Range("L115").Value = 3500
$ws.Range("M115").Value = 3500
$ws.Range("P115").Value = 3500

# Row 116
$ws.Range("D116").Value = 44188
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 3000
$ws.Range("P116").Value = 3000

# Row 117
$ws.Range("D117").Value = 44188
$ws.Range("J117").Value = 3000
$ws.Range("K117").Value = 2500
$ws.Range("L117").Value = 2500
$ws.Range("M117").Value = 2500
$ws.Range("P117").Value = 2500

# Row 118
$ws.Range("D118").Value = 44224
$ws.Range("K118").Value = 2300
$ws.Range("L118").Value = 2300
$ws.Range("M118").Value = 2300
$ws.Range("P118").Value = 2300

# Row 119
$ws.Range("D119").Value = 44224
$ws.Range("J119").Value = 6000
$ws.Range("K119").Value = 1800
$ws.Range("L119").Value = 1800
$ws.Range("M119").Value = 1800
$ws.Range("P119").Value = 1800

# Row 120
$ws.Range("D120").Value = 44224
$ws.Range("J120").Value = 4000
$ws.Range("K120").Value = 1300
$ws.Range("L120").Value = 1300
$ws.Range("M120").Value = 1300
$ws.Range("P120").Value = 1300

# Row 121
$ws.Range("D121").Value = 44195

# Row 122
$ws.Range("D122").Value = 44195
$ws.Range("J122").Value = 3000

# Row 123
$ws.Range("D123").Value = 44195

# --- Append new rows 124-126 (full new records) ---
# Row 124
$ws.Range("D124").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A124").Value = 5
$ws.Range("B124").Value = 'Macroferia Regional de Talca'
$ws.Range("C124").Value = 'Maule'
$ws.Range("D124").Value = 44194
$ws.Range("E124").Value = 7
$ws.Range("F124").Value = 100112028
$ws.Range("G124").Value = 'Sandia'
$ws.Range("H124").Value = 'Sin especificar'
$ws.Range("I124").Value = 'Extra'
$ws.Range("J124").Value = 2000
$ws.Range("K124").Value = 3000
$ws.Range("L124").Value = 3000
$ws.Range("M124").Value = 3000
$ws.Range("N124").Value = '$/unidad'
$ws.Range("O124").Value = 'Región del Maule'
$ws.Range("P124").Value = 3000
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = 'Hortaliza'

# Row 125
$ws.Range("D125").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A125").Value = 5
$ws.Range("B125").Value = 'Macroferia Regional de Talca'
$ws.Range("C125").Value = 'Maule'
$ws.Range("D125").Value = 44194
$ws.Range("E125").Value = 7
$ws.Range("F125").Value = 100112028
$ws.Range("G125").Value = 'Sandia'
$ws.Range("H125").Value = 'Sin especificar'
$ws.Range("I125").Value = 'Primera'
$ws.Range("J125").Value = 4000
$ws.Range("K125").Value = 2500
$ws.Range("L125").Value = 2500
$ws.Range("M125").Value = 2500
$ws.Range("N125").Value = '$/unidad'
$ws.Range("O125").Value = 'Región del Maule'
$ws.Range("P125").Value = 2500
$ws.Range("Q125").Value = 1
$ws.Range("R125").Value = 'Hortaliza'

# Row 126
$ws.Range("D126").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A126").Value = 5
$ws.Range("B126").Value = 'Macroferia Regional de Talca'
$ws.Range("C126").Value = 'Maule'
$ws.Range("D126").Value = 44194
$ws.Range("E126").Value = 7
$ws.Range("F126").Value = 100112028
$ws.Range("G126").Value = 'Sandia'
$ws.Range("H126").Value = 'Sin especificar'
$ws.Range("I126").Value = 'Segunda'
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 2000
$ws.Range("L126").Value = 2000
$ws.Range("M126").Value = 2000
$ws.Range("N126").Value = '$/unidad'
$ws.Range("O126").Value = 'Región del Maule'
$ws.Range("P126").Value = 2000
$ws.Range("Q126").Value = 1
$ws.Range("R126").Value = 'Hortaliza'

